# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# The workbook tracks, per-quarter, which mutual funds hold the stock and
# what share of the stock they own. A new quarter sheet "2022-Q1" is
# inserted right before the "总计" (total/summary) sheet, and the "总计"
# sheet gets a new leading row summarizing the 2022-Q1 numbers.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Detailed per-fund holdings for the new "2022-Q1" sheet (columns B..H).
# Column A is a simple 0-based row counter, filled in separately below.
$q1Data = @(
    @("003634","嘉实农业产业股票","28.18","94.29","7.42","2.0910","7"),
    @("161810","银华内需精选混合(LOF)","25.59","94.71","7.98","2.0421","3"),
    @("003751","万家瑞隆混合","27.84","86.40","5.94","1.6537","4"),
    @("009394","银华同力精选混合","20.03","94.68","7.75","1.5523","3"),
    @("009199","万家价值优势一年持有期混合","13.70","89.53","8.61","1.1796","5"),
    @("161838","银华创业板两年定期开放混合","10.44","95.40","8.29","0.8655","2"),
    @("005094","万家臻选混合","13.43","73.12","4.77","0.6406","9"),
    @("005106","银华农业产业股票","13.24","93.41","4.25","0.5627","10"),
    @("161912","万家社会责任18个月定期开放混合（LOF）A","13.56","88.11","4.11","0.5573","9"),
    @("180020","银华成长先锋混合","3.05","79.81","6.60","0.2013","3"),
    @("900008","中信证券稳健回报混合A","6.36","88.19","2.99","0.1902","7"),
    @("900078","中信证券稳健回报混合C","1.66","88.19","2.99","0.0496","7"),
    @("900027","中信证券信远一年持有期混合型集合资产管理计划A","0.71","75.94","3.32","0.0236","2"),
    @("161913","万家社会责任18个月定期开放混合（LOF）C","0.44","88.11","4.11","0.0181","9"),
    @("900087","中信证券信远一年持有期混合型集合资产管理计划C","0.02","75.94","3.32","0.0007","2"),
    @("900077","中信证券信远一年持有期混合型集合资产管理计划B","0.01","75.94","3.32","0.0003","2")
)

# Summary rows for the "总计" sheet (columns B..D): quarter label, number of
# holding funds, total held market value (亿元). The first row (2022-Q1) is
# the newly added one; the rest already existed and simply shift down.
$totalData = @(
    @("2022-Q1","16","11.63"),
    @("2021-Q4","9","10.25"),
    @("2021-Q3","12","10.58"),
    @("2021-Q2","12","10.38"),
    @("2021-Q1","10","9.93"),
    @("2020-Q4","17","13.33")
)

# ---------------------------------------------------------------------------
# Re-create the "总计" sheet and insert a fresh "2022-Q1" sheet right before
# it, so the sheetId/relationship-id numbering matches what a natural
# "insert sheet, then keep total sheet last" edit would produce.
# ---------------------------------------------------------------------------

$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQ1 = $wb.Worksheets.Add($null, $lastSheet)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# ---------------------------------------------------------------------------
# Style helpers: the rest of the workbook formats header-row and first
# column cells with a bold/centered/bordered style. Reuse that exact style
# (instead of re-building it) by copying formats from the "2021-Q4" sheet,
# which already uses it throughout.
# ---------------------------------------------------------------------------

$styleSrcSheet = $wb.Worksheets.Item("2021-Q4")

# Match page margins used across all the other sheets in the workbook.
$wsQ1.PageSetup.LeftMargin = 0.75 * 72
$wsQ1.PageSetup.RightMargin = 0.75 * 72
$wsQ1.PageSetup.TopMargin = 1 * 72
$wsQ1.PageSetup.BottomMargin = 1 * 72
$wsQ1.PageSetup.HeaderMargin = 0.5 * 72
$wsQ1.PageSetup.FooterMargin = 0.5 * 72
$wsQ1.Outline.SummaryRow = 1
$wsQ1.Outline.SummaryColumn = 1

$wsTotal.PageSetup.LeftMargin = 0.75 * 72
$wsTotal.PageSetup.RightMargin = 0.75 * 72
$wsTotal.PageSetup.TopMargin = 1 * 72
$wsTotal.PageSetup.BottomMargin = 1 * 72
$wsTotal.PageSetup.HeaderMargin = 0.5 * 72
$wsTotal.PageSetup.FooterMargin = 0.5 * 72
$wsTotal.Outline.SummaryRow = 1
$wsTotal.Outline.SummaryColumn = 1

# ---------------------------------------------------------------------------
# Populate "2022-Q1"
# ---------------------------------------------------------------------------

$q1HeaderSrc = $styleSrcSheet.Range("B1:H1")
$q1HeaderDst = $wsQ1.Range("B1:H1")
$q1HeaderSrc.Copy()
$q1HeaderDst.PasteSpecial(-4122)

$wsQ1.Cells.Item(1,2).Value = "基金代码"
$wsQ1.Cells.Item(1,3).Value = "基金名称"
$wsQ1.Cells.Item(1,4).Value = "基金规模"
$wsQ1.Cells.Item(1,5).Value = "股票总仓位"
$wsQ1.Cells.Item(1,6).Value = "仓位占比"
$wsQ1.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ1.Cells.Item(1,8).Value = "仓位排名"

$q1RowCount = $q1Data.Count

$q1ColASrc = $styleSrcSheet.Range("A2:A10")
$q1ColADst = $wsQ1.Range("A2:A$(1 + $q1RowCount)")
$q1ColASrc.Copy()
$q1ColADst.PasteSpecial(-4122)

# Columns B..G hold values that look numeric ("003634", "28.18", ...) but
# must stay text so leading zeros and fixed decimal digits are preserved,
# exactly like the other per-quarter sheets in this workbook. Force text
# format before writing, then drop back to the Normal style (no explicit
# number format) so the saved cells have no extra style index.
$q1TextRange = $wsQ1.Range("B2:G$(1 + $q1RowCount)")
$q1TextRange.NumberFormat = "@"

for ($i = 0; $i -lt $q1RowCount; $i++) {
    $r = 2 + $i
    $row = $q1Data[$i]

    $wsQ1.Cells.Item($r, 1).Value = $i

    $wsQ1.Cells.Item($r, 2).Value = $row[0]
    $wsQ1.Cells.Item($r, 3).Value = $row[1]
    $wsQ1.Cells.Item($r, 4).Value = $row[2]
    $wsQ1.Cells.Item($r, 5).Value = $row[3]
    $wsQ1.Cells.Item($r, 6).Value = $row[4]
    $wsQ1.Cells.Item($r, 7).Value = $row[5]
    $wsQ1.Cells.Item($r, 8).Value = [int]$row[6]
}

$q1TextRange.Style = "Normal"

# ---------------------------------------------------------------------------
# Populate "总计"
# ---------------------------------------------------------------------------

$totalHeaderSrc = $styleSrcSheet.Range("B1:D1")
$totalHeaderDst = $wsTotal.Range("B1:D1")
$totalHeaderSrc.Copy()
$totalHeaderDst.PasteSpecial(-4122)

$wsTotal.Cells.Item(1,2).Value = "日期"
$wsTotal.Cells.Item(1,3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalRowCount = $totalData.Count

$totalColASrc = $styleSrcSheet.Range("A2:A10")
$totalColADst = $wsTotal.Range("A2:A$(1 + $totalRowCount)")
$totalColASrc.Copy()
$totalColADst.PasteSpecial(-4122)

for ($i = 0; $i -lt $totalRowCount; $i++) {
    $r = 2 + $i
    $row = $totalData[$i]

    $wsTotal.Cells.Item($r, 1).Value = $i
    $wsTotal.Cells.Item($r, 2).Value = $row[0]
    $wsTotal.Cells.Item($r, 3).Value = [int]$row[1]
    $wsTotal.Cells.Item($r, 4).Value = [double]$row[2]
}

# ---------------------------------------------------------------------------
# Restore the originally-active sheet ("2020-Q4") as the selected tab.
# ---------------------------------------------------------------------------

$wb.Worksheets.Item(1).Activate()
